# Updated symbol list on Wed Feb  1 10:48:04 UTC 2023 with GitHub Actions
# Refresh Price (col D) and Volume(1h) (col E) figures for the crypto
# table. Values are stored as literal text in the sheet (e.g. "307.11",
# "-1.03%"), so each cell is first forced to Text format ("@") before the
# new value is written — otherwise Excel would auto-convert the numeric-
# looking strings into real numbers/percentages and lose the exact
# formatting (trailing zeros, literal "%" characters, etc.).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.11"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.03%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.96"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.56%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.116"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.41%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.217"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.02%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.875"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.19%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.990"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.11%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9229"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.40%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1101"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-7.88%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1892"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.30%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08857"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.40%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-2.04%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09586"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001385"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.40%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005888"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.46%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.427"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.51%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.409"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.45%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.47%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.237"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "18.48%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1287"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.42%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-6.11%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04337"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.89%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001195"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.39%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004247"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.05%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-2.32%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-98.11%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02145"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.97%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05003"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.43%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007529"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.49%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1348"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.26%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008686"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-11.66%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.001993"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-8.29%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-8.65%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006546"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.59%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003384"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "15.20%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-16.60%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002004"
